# Add a new "2022-Q4" quarter sheet to the workbook, insert its summary
# row into the "总计" (total) sheet, and shift the existing quarter rows
# down to make room for it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "总计" summary sheet: insert the new 2022-Q4 row and push the
#    existing 2022-Q3 / 2022-Q2 / 2022-Q1 rows down by one.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Clone the formatted style of row 4 into the brand-new row 5 before we
# start overwriting values, so the new last row keeps the same look
# (bordered / bold / centered) as column A of the other data rows.
$summary.Range("A4").Copy()
$summary.Range("A5").PasteSpecial(-4122)

# Shift data down one row at a time, starting from the bottom so we
# never clobber a value before it has been copied onward.
$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2022-Q1"
$summary.Range("C5").Value = 1
$summary.Range("D5").Value = 0.9

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2022-Q2"
$summary.Range("C4").Value = 3
$summary.Range("D4").Value = 1.79

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 5
$summary.Range("D3").Value = 1.33

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.06

# ---------------------------------------------------------------------
# 2. Create the new "2022-Q4" worksheet, positioned right after "总计"
#    (i.e. right before the existing "2022-Q3" sheet).
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($wb.Worksheets.Item(2))
$q4.Name = "2022-Q4"

# Re-fetch the "2022-Q3" sheet by name now that the sheet collection has
# shifted (the worksheet handle captured before Add() above tracks the
# *position* it was grabbed from, not the sheet itself, once a sheet is
# inserted ahead of it).
$existingQ3 = $wb.Worksheets.Item("2022-Q3")

# Borrow the header/row styling from the "2022-Q3" sheet so the new
# sheet matches the look of the other quarterly breakdown sheets.
$existingQ3.Range("B1:H1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$existingQ3.Range("A2:A3").Copy()
$q4.Range("A2:A3").PasteSpecial(-4122)

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

$q4.Range("A2").Value = 0
$q4.Range("B2").NumberFormat = "@"
$q4.Range("B2").Value = "014232"
$q4.Range("C2").Value = "博时专精特新主题混合A"
$q4.Range("D2").NumberFormat = "@"
$q4.Range("D2").Value = "3.00"
$q4.Range("E2").NumberFormat = "@"
$q4.Range("E2").Value = "80.89"
$q4.Range("F2").NumberFormat = "@"
$q4.Range("F2").Value = "1.13"
$q4.Range("G2").NumberFormat = "@"
$q4.Range("G2").Value = "0.0339"
$q4.Range("H2").Value = 8

$q4.Range("A3").Value = 1
$q4.Range("B3").NumberFormat = "@"
$q4.Range("B3").Value = "014233"
$q4.Range("C3").Value = "博时专精特新主题混合C"
$q4.Range("D3").NumberFormat = "@"
$q4.Range("D3").Value = "2.58"
$q4.Range("E3").NumberFormat = "@"
$q4.Range("E3").Value = "80.89"
$q4.Range("F3").NumberFormat = "@"
$q4.Range("F3").Value = "1.13"
$q4.Range("G3").NumberFormat = "@"
$q4.Range("G3").Value = "0.0292"
$q4.Range("H3").Value = 8

# Match the page margins used by the other worksheets (0.75in / 1in /
# 0.5in left-right / top-bottom / header-footer).
$q4.PageSetup.LeftMargin = 54
$q4.PageSetup.RightMargin = 54
$q4.PageSetup.TopMargin = 72
$q4.PageSetup.BottomMargin = 72
$q4.PageSetup.HeaderMargin = 36
$q4.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------
# 3. Restore the original active sheet ("2022-Q1", now the last tab)
#    since adding a worksheet moves the selection to the new sheet.
# ---------------------------------------------------------------------
$wb.Worksheets.Item($wb.Worksheets.Count).Activate()
